# Update loading_percent values for rows 2-25, columns B:O (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 24,14
$data[0,0] = 9.962056375687959
$data[0,1] = 6.01085826928519
$data[0,2] = 14.81917940431809
$data[0,3] = 16.06956077238866
$data[0,4] = 0
$data[0,5] = 3.697537859662062
$data[0,6] = 0
$data[0,7] = 27.09372476323565
$data[0,8] = 9.371417271002624
$data[0,9] = 10.67763446192085
$data[0,10] = 0
$data[0,11] = 17.02962569999848
$data[0,12] = 0
$data[0,13] = 30.22109951554191
$data[1,0] = 9.745444507363054
$data[1,1] = 5.879221636905421
$data[1,2] = 14.79467827794214
$data[1,3] = 16.06953503015991
$data[1,4] = 0
$data[1,5] = 3.699649961894595
$data[1,6] = 0
$data[1,7] = 27.17531298897445
$data[1,8] = 9.389179777339313
$data[1,9] = 10.54606222751534
$data[1,10] = 0
$data[1,11] = 16.98456554604114
$data[1,12] = 0
$data[1,13] = 30.29184192596547
$data[2,0] = 9.611972179904537
$data[2,1] = 5.798124349489171
$data[2,2] = 14.78259053086428
$data[2,3] = 16.07243075518878
$data[2,4] = 0
$data[2,5] = 3.701015726755772
$data[2,6] = 0
$data[2,7] = 27.22993834682391
$data[2,8] = 9.400996131946386
$data[2,9] = 10.46643987151869
$data[2,10] = 0
$data[2,11] = 16.95975818794103
$data[2,12] = 0
$data[2,13] = 30.3408254635133
$data[3,0] = 9.557543400939377
$data[3,1] = 5.765056959105439
$data[3,2] = 14.77841182059869
$data[3,3] = 16.07434394779656
$data[3,4] = 0
$data[3,5] = 3.701589674472422
$data[3,6] = 0
$data[3,7] = 27.25333652746239
$data[3,8] = 9.40604057292688
$data[3,9] = 10.43432239838666
$data[3,10] = 0
$data[3,11] = 16.95037550914745
$data[3,12] = 0
$data[3,13] = 30.36217942344215
$data[4,0] = 9.54850567505204
$data[4,1] = 5.759566413351813
$data[4,2] = 14.777763175015
$data[4,3] = 16.07470593706982
$data[4,4] = 0
$data[4,5] = 3.701686029846857
$data[4,6] = 0
$data[4,7] = 27.25729046959277
$data[4,8] = 9.406892050046732
$data[4,9] = 10.42901026683958
$data[4,10] = 0
$data[4,11] = 16.9488616067689
$data[4,12] = 0
$data[4,13] = 30.36580928564685
$data[5,0] = 9.611238178474713
$data[5,1] = 5.79767840494587
$data[5,2] = 14.78253114553202
$data[5,3] = 16.072453587501
$data[5,4] = 0
$data[5,5] = 3.701023396740714
$data[5,6] = 0
$data[5,7] = 27.23024929712745
$data[5,8] = 9.401063234665195
$data[5,9] = 10.46600534267925
$data[5,10] = 0
$data[5,11] = 16.9596286988786
$data[5,12] = 0
$data[5,13] = 30.34110781422761
$data[6,0] = 9.887514292650849
$data[6,1] = 5.96555519022188
$data[6,2] = 14.81012035248282
$data[6,3] = 16.06894831782442
$data[6,4] = 0
$data[6,5] = 3.698251840510162
$data[6,6] = 0
$data[6,7] = 27.12091551980344
$data[6,8] = 9.377353088380486
$data[6,9] = 10.63204902451388
$data[6,10] = 0
$data[6,11] = 17.01350004631741
$data[6,12] = 0
$data[6,13] = 30.24433887207057
$data[7,0] = 10.42212908126817
$data[7,1] = 6.290542300773391
$data[7,2] = 14.88748550432578
$data[7,3] = 16.08511933197894
$data[7,4] = 0
$data[7,5] = 3.693361185351387
$data[7,6] = 0
$data[7,7] = 26.9425061297699
$data[7,8] = 9.338064322379667
$data[7,9] = 10.9652622899662
$data[7,10] = 0
$data[7,11] = 17.14148078304184
$data[7,12] = 0
$data[7,13] = 30.09868384498784
$data[8,0] = 10.80619758198532
$data[8,1] = 6.524116344752352
$data[8,2] = 14.95821513718033
$data[8,3] = 16.11095565988959
$data[8,4] = 0
$data[8,5] = 3.690096327648641
$data[8,6] = 0
$data[8,7] = 26.83344084485036
$data[8,8] = 9.313572613289988
$data[8,9] = 11.21251963769148
$data[8,10] = 0
$data[8,11] = 17.24861861201837
$data[8,12] = 0
$data[8,13] = 30.01867997884025
$data[9,0] = 10.97818459405884
$data[9,1] = 6.628738317533384
$data[9,2] = 14.99333379108702
$data[9,3] = 16.1257128955559
$data[9,4] = 0
$data[9,5] = 3.688681593089384
$data[9,6] = 0
$data[9,7] = 26.78861834516897
$data[9,8] = 9.303376236481704
$data[9,9] = 11.32507454009611
$data[9,10] = 0
$data[9,11] = 17.3000848723137
$data[9,12] = 0
$data[9,13] = 29.98817093004528
$data[10,0] = 11.04285479451945
$data[10,1] = 6.668082286579974
$data[10,2] = 15.00704838189694
$data[10,3] = 16.13173019431465
$data[10,4] = 0
$data[10,5] = 3.688155945534826
$data[10,6] = 0
$data[10,7] = 26.77233547710563
$data[10,8] = 9.299650696328012
$data[10,9] = 11.3676694510473
$data[10,10] = 0
$data[10,11] = 17.31995512568908
$data[10,12] = 0
$data[10,13] = 29.97746569167811
$data[11,0] = 11.02894825172026
$data[11,1] = 6.659621648544162
$data[11,2] = 15.00407632381441
$data[11,3] = 16.13041522959049
$data[11,4] = 0
$data[11,5] = 3.688268705563942
$data[11,6] = 0
$data[11,7] = 26.77581155505607
$data[11,8] = 9.300447031506152
$data[11,9] = 11.35849768207126
$data[11,10] = 0
$data[11,11] = 17.31565894684505
$data[11,12] = 0
$data[11,13] = 29.97973352347005
$data[12,0] = 10.98351457853525
$data[12,1] = 6.631980879133867
$data[12,2] = 14.99445380702943
$data[12,3] = 16.12619936154305
$data[12,4] = 0
$data[12,5] = 3.688638145987582
$data[12,6] = 0
$data[12,7] = 26.78726489861951
$data[12,8] = 9.303067017943189
$data[12,9] = 11.32857960933883
$data[12,10] = 0
$data[12,11] = 17.30171203872893
$data[12,12] = 0
$data[12,13] = 29.98727320250016
$data[13,0] = 10.95562369340186
$data[13,1] = 6.61501328983074
$data[13,2] = 14.9886136756185
$data[13,3] = 16.12367280291924
$data[13,4] = 0
$data[13,5] = 3.688865750340934
$data[13,6] = 0
$data[13,7] = 26.79437035812965
$data[13,8] = 9.30468948671481
$data[13,9] = 11.3102492678544
$data[13,10] = 0
$data[13,11] = 17.29321844117723
$data[13,12] = 0
$data[13,13] = 29.99200193614646
$data[14,0] = 10.79489725369284
$data[14,1] = 6.517242781238017
$data[14,2] = 14.95597866424638
$data[14,3] = 16.11005144777243
$data[14,4] = 0
$data[14,5] = 3.69019019744579
$data[14,6] = 0
$data[14,7] = 26.83646664350526
$data[14,8] = 9.314257962301815
$data[14,9] = 11.20516216383236
$data[14,10] = 0
$data[14,11] = 17.24530912649174
$data[14,12] = 0
$data[14,13] = 30.02079238235334
$data[15,0] = 10.69555028524088
$data[15,1] = 6.456816841599686
$data[15,2] = 14.93670696685163
$data[15,3] = 16.10246264321101
$data[15,4] = 0
$data[15,5] = 3.6910207142551
$data[15,6] = 0
$data[15,7] = 26.8635196506286
$data[15,8] = 9.320369751797205
$data[15,9] = 11.1406867545874
$data[15,10] = 0
$data[15,11] = 17.21660933810452
$data[15,12] = 0
$data[15,13] = 30.03996291468006
$data[16,0] = 10.63815542103273
$data[16,1] = 6.421910013337147
$data[16,2] = 14.92589991609173
$data[16,3] = 16.09838072743867
$data[16,4] = 0
$data[16,5] = 3.691505041143557
$data[16,6] = 0
$data[16,7] = 26.87953077584143
$data[16,8] = 9.32397405885837
$data[16,9] = 11.10361162067041
$data[16,10] = 0
$data[16,11] = 17.20035952882396
$data[16,12] = 0
$data[16,13] = 30.05154314935338
$data[17,0] = 10.61868109632444
$data[17,1] = 6.410066398568748
$data[17,2] = 14.92228871335931
$data[17,3] = 16.09704734412765
$data[17,4] = 0
$data[17,5] = 3.691670167139446
$data[17,6] = 0
$data[17,7] = 26.8850292764117
$data[17,8] = 9.325209704261818
$data[17,9] = 11.09106141117724
$data[17,10] = 0
$data[17,11] = 17.19490219256468
$data[17,12] = 0
$data[17,13] = 30.0555590900997
$data[18,0] = 10.70615262339949
$data[18,1] = 6.46326525212103
$data[18,2] = 14.93872980110987
$data[18,3] = 16.10324121891312
$data[18,4] = 0
$data[18,5] = 3.690931617936232
$data[18,6] = 0
$data[18,7] = 26.8605931299605
$data[18,8] = 9.319709935653774
$data[18,9] = 11.14754957658956
$data[18,10] = 0
$data[18,11] = 17.21963790426658
$data[18,12] = 0
$data[18,13] = 30.03786484491168
$data[19,0] = 10.99687246575376
$data[19,1] = 6.640107382925077
$data[19,2] = 14.9972689466825
$data[19,3] = 16.12742604493901
$data[19,4] = 0
$data[19,5] = 3.688529359170952
$data[19,6] = 0
$data[19,7] = 26.78388202485544
$data[19,8] = 9.302293786389379
$data[19,9] = 11.3373683100788
$data[19,10] = 0
$data[19,11] = 17.30579833333771
$data[19,12] = 0
$data[19,13] = 29.98503559314421
$data[20,0] = 11.18418061418997
$data[20,1] = 6.75406989126471
$data[20,2] = 15.03794817761092
$data[20,3] = 16.14573147307901
$data[20,4] = 0
$data[20,5] = 3.687018084723669
$data[20,6] = 0
$data[20,7] = 26.73777190560042
$data[20,8] = 9.291701610282912
$data[20,9] = 11.46125031983891
$data[20,10] = 0
$data[20,11] = 17.36432489455752
$data[20,12] = 0
$data[20,13] = 29.95545128889136
$data[21,0] = 11.0844776583658
$data[21,1] = 6.693405964967203
$data[21,2] = 15.01601799069206
$data[21,3] = 16.13573390129283
$data[21,4] = 0
$data[21,5] = 3.687819322372017
$data[21,6] = 0
$data[21,7] = 26.76201300925707
$data[21,8] = 9.297282638051339
$data[21,9] = 11.39516058637157
$data[21,10] = 0
$data[21,11] = 17.33288928641992
$data[21,12] = 0
$data[21,13] = 29.9707882594362
$data[22,0] = 10.70136017022661
$data[22,1] = 6.460350443332019
$data[22,2] = 14.93781442837605
$data[22,3] = 16.1028883493716
$data[22,4] = 0
$data[22,5] = 3.690971877035115
$data[22,6] = 0
$data[22,7] = 26.8619147834021
$data[22,8] = 9.320007956427968
$data[22,9] = 11.14444691476819
$data[22,10] = 0
$data[22,11] = 17.21826790967238
$data[22,12] = 0
$data[22,13] = 30.03881164157738
$data[23,0] = 10.2787315188736
$data[23,1] = 6.203355683073135
$data[23,2] = 14.86409359450581
$data[23,3] = 16.07828581729013
$data[23,4] = 0
$data[23,5] = 3.694626328170079
$data[23,6] = 0
$data[23,7] = 26.98690993995056
$data[23,8] = 9.347923473605995
$data[23,9] = 10.8745319882276
$data[23,10] = 0
$data[23,11] = 17.10451711644402
$data[23,12] = 0
$data[23,13] = 30.09868384498784
$ws.Range("B2:O25").Value2 = $data
